$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 7-12: update tp (C), fp (D), tn (E), fn (F) values and the
# recomputed precision (G), recall (H), fscore (I) columns.
$data = @(
    @{ Row = 7;  C = 24; D = 133; E = 0; F = 0; G = 0.1528662420382166; H = 1; I = 0.2651933701657458 },
    @{ Row = 8;  C = 24; D = 133; E = 0; F = 0; G = 0.1528662420382166; H = 1; I = 0.2651933701657458 },
    @{ Row = 9;  C = 24; D = 130; E = 0; F = 0; G = 0.1558441558441558; H = 1; I = 0.2696629213483146 },
    @{ Row = 10; C = 24; D = 128; E = 0; F = 0; G = 0.1578947368421053; H = 1; I = 0.2727272727272727 },
    @{ Row = 11; C = 24; D = 118; E = 0; F = 0; G = 0.1690140845070423; H = 1; I = 0.2891566265060241 },
    @{ Row = 12; C = 23; D = 106; E = 0; F = 0; G = 0.1782945736434109; H = 1; I = 0.3026315789473684 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("C$r").Value = $entry.C
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("E$r").Value = $entry.E
    $ws.Range("F$r").Value = $entry.F
    $ws.Range("G$r").Value = $entry.G
    $ws.Range("H$r").Value = $entry.H
    $ws.Range("I$r").Value = $entry.I
}

$wb.Save()
